# Desalin test data changed
# Update existing row 2 (Amount + DateofTransfer columns) and append new
# rows 3-6 of test data to the LocalFT sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: change DateofTransfer (H2) ---
# H2 switches from a date-formatted cell to a plain text cell holding "1111".
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "'1111"

# --- New rows 3-6 ---

# Column A: Databinding key
$ws.Range("A3").Value = "TestData_002"
$ws.Range("A4").Value = "TestData_003"
$ws.Range("A5").Value = "TestData_004"
$ws.Range("A6").Value = "TestData_005"

# Column B: ReceiverBankName
$ws.Range("B3").Value = "SBI"
$ws.Range("B4").Value = "SCB"
$ws.Range("B5").Value = "ICICI"
$ws.Range("B6").Value = "IOB"

# Column C: ReceiverName
$ws.Range("C3").Value = "kevin"
$ws.Range("C4").Value = "mohan"
$ws.Range("C5").Value = "sambu"
$ws.Range("C6").Value = "gayathri"

# Column E: SWIFTmsg
$ws.Range("E3").Value = "MT104"
$ws.Range("E4").Value = "MT105"
$ws.Range("E5").Value = "MT106"
$ws.Range("E6").Value = "MT107"

# Column F: Amount (quote-prefixed text). F2 changes 60000 -> 6, and new
# rows 3-6 follow the same quote-prefixed text style.
$ws.Range("F2").Value = "'6"
$ws.Range("F3").Value = "'7"
$ws.Range("F4").Value = "'8"
$ws.Range("F5").Value = "'9"
$ws.Range("F6").Value = "'10"

# Column G: TransferType
$ws.Range("G3").Value = "Local Transfer"
$ws.Range("G4").Value = "Local Transfer"
$ws.Range("G5").Value = "Local Transfer"
$ws.Range("G6").Value = "Local Transfer"

# Column H: DateofTransfer (now plain text, same style as H2)
$ws.Range("H3").Value = "'1112"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H4").Value = "'1113"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H5").Value = "'1114"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H6").Value = "'1115"
$ws.Range("H6").NumberFormat = "@"

# Column I: TransferDescription
$ws.Range("I3").Value = "fund transfer"
$ws.Range("I4").Value = "fund transfer"
$ws.Range("I5").Value = "fund transfer"
$ws.Range("I6").Value = "fund transfer"

# Column D: ReceiverAccNum (plain numeric, no special style)
$ws.Range("D3").Value = 25632523
$ws.Range("D4").Value = 23365478
$ws.Range("D5").Value = 23657489
$ws.Range("D6").Value = 1236548

# Widen column H to fit the new values and drop its autofit flag.
$ws.Columns("H").ColumnWidth = 25.5

# Selection now spans the freshly populated rows.
$ws.Range("A2:A6").Select() | Out-Null
